$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 368; $r++) {
    $ws.Cells.Item($r, 3).Value = 46061
}
